$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Re-run SGNN annotation: update DAMSLTag (col I) and DialogAct (col J)
# values for the affected rows following transcript clean-up.

$ws.Range("I2").Value = 'sd'
$ws.Range("J2").Value = 'Statement-non-opinion'
$ws.Range("I3").Value = 'aa'
$ws.Range("J3").Value = 'Agree/Accept'
$ws.Range("I10").Value = 'ba'
$ws.Range("J10").Value = 'Appreciation'
$ws.Range("I12").Value = 'b'
$ws.Range("J12").Value = 'Acknowledge (Backchannel)'
$ws.Range("I18").Value = 'sv'
$ws.Range("J18").Value = 'Statement-opinion'
$ws.Range("I19").Value = 'sd'
$ws.Range("J19").Value = 'Statement-non-opinion'
$ws.Range("I21").Value = 'b'
$ws.Range("J21").Value = 'Acknowledge (Backchannel)'
$ws.Range("I27").Value = 'sv'
$ws.Range("J27").Value = 'Statement-opinion'
$ws.Range("I55").Value = 'b'
$ws.Range("J55").Value = 'Acknowledge (Backchannel)'
$ws.Range("I58").Value = 'sv'
$ws.Range("J58").Value = 'Statement-opinion'
$ws.Range("I71").Value = 'sd'
$ws.Range("J71").Value = 'Statement-non-opinion'
$ws.Range("I74").Value = 'b'
$ws.Range("J74").Value = 'Acknowledge (Backchannel)'
$ws.Range("I75").Value = 'sv'
$ws.Range("J75").Value = 'Statement-opinion'
$ws.Range("I89").Value = '%'
$ws.Range("J89").Value = 'Uninterpretable'
$ws.Range("I106").Value = 'sv'
$ws.Range("J106").Value = 'Statement-opinion'
$ws.Range("I114").Value = 'sd'
$ws.Range("J114").Value = 'Statement-non-opinion'
$ws.Range("I115").Value = 'sv'
$ws.Range("J115").Value = 'Statement-opinion'
$ws.Range("I116").Value = 'aa'
$ws.Range("J116").Value = 'Agree/Accept'
$ws.Range("I123").Value = 'sd'
$ws.Range("J123").Value = 'Statement-non-opinion'
$ws.Range("I138").Value = 'sv'
$ws.Range("J138").Value = 'Statement-opinion'
$ws.Range("I159").Value = 'ba'
$ws.Range("J159").Value = 'Appreciation'
$ws.Range("I164").Value = '%'
$ws.Range("J164").Value = 'Uninterpretable'
$ws.Range("I166").Value = 'sv'
$ws.Range("J166").Value = 'Statement-opinion'
$ws.Range("I170").Value = 'aa'
$ws.Range("J170").Value = 'Agree/Accept'
$ws.Range("I174").Value = 'aa'
$ws.Range("J174").Value = 'Agree/Accept'
$ws.Range("I181").Value = 'aa'
$ws.Range("J181").Value = 'Agree/Accept'
$ws.Range("I186").Value = 'aa'
$ws.Range("J186").Value = 'Agree/Accept'
$ws.Range("I188").Value = 'sd'
$ws.Range("J188").Value = 'Statement-non-opinion'
$ws.Range("I190").Value = 'qy'
$ws.Range("J190").Value = 'Yes-No-Question'
$ws.Range("I202").Value = 'sd'
$ws.Range("J202").Value = 'Statement-non-opinion'
$ws.Range("I209").Value = 'b'
$ws.Range("J209").Value = 'Acknowledge (Backchannel)'
$ws.Range("I231").Value = 'b'
$ws.Range("J231").Value = 'Acknowledge (Backchannel)'
$ws.Range("I284").Value = 'sd'
$ws.Range("J284").Value = 'Statement-non-opinion'
$ws.Range("I292").Value = 'sv'
$ws.Range("J292").Value = 'Statement-opinion'
$ws.Range("I301").Value = 'sv'
$ws.Range("J301").Value = 'Statement-opinion'
$ws.Range("I304").Value = 'sv'
$ws.Range("J304").Value = 'Statement-opinion'
$ws.Range("I307").Value = 'aa'
$ws.Range("J307").Value = 'Agree/Accept'
$ws.Range("I308").Value = 'aa'
$ws.Range("J308").Value = 'Agree/Accept'
$ws.Range("I309").Value = 'aa'
$ws.Range("J309").Value = 'Agree/Accept'
$ws.Range("I312").Value = 'b'
$ws.Range("J312").Value = 'Acknowledge (Backchannel)'
$ws.Range("I313").Value = 'sv'
$ws.Range("J313").Value = 'Statement-opinion'
$ws.Range("I329").Value = 'aa'
$ws.Range("J329").Value = 'Agree/Accept'
$ws.Range("I337").Value = 'sv'
$ws.Range("J337").Value = 'Statement-opinion'
$ws.Range("I346").Value = 'sv'
$ws.Range("J346").Value = 'Statement-opinion'
$ws.Range("I348").Value = 'b'
$ws.Range("J348").Value = 'Acknowledge (Backchannel)'
$ws.Range("I350").Value = 'aa'
$ws.Range("J350").Value = 'Agree/Accept'
$ws.Range("I352").Value = '%'
$ws.Range("J352").Value = 'Uninterpretable'
$ws.Range("I360").Value = 'ba'
$ws.Range("J360").Value = 'Appreciation'
$ws.Range("I362").Value = '%'
$ws.Range("J362").Value = 'Uninterpretable'
$ws.Range("I364").Value = 'b'
$ws.Range("J364").Value = 'Acknowledge (Backchannel)'
$ws.Range("I370").Value = 'sd'
$ws.Range("J370").Value = 'Statement-non-opinion'
$ws.Range("I384").Value = 'sv'
$ws.Range("J384").Value = 'Statement-opinion'
$ws.Range("I390").Value = 'sd'
$ws.Range("J390").Value = 'Statement-non-opinion'
$ws.Range("I414").Value = 'sd'
$ws.Range("J414").Value = 'Statement-non-opinion'
$ws.Range("I417").Value = 'sd'
$ws.Range("J417").Value = 'Statement-non-opinion'
$ws.Range("I429").Value = '%'
$ws.Range("J429").Value = 'Uninterpretable'
$ws.Range("I440").Value = 'b'
$ws.Range("J440").Value = 'Acknowledge (Backchannel)'
$ws.Range("I455").Value = 'aa'
$ws.Range("J455").Value = 'Agree/Accept'
$ws.Range("I456").Value = 'aa'
$ws.Range("J456").Value = 'Agree/Accept'
$ws.Range("I460").Value = 'sd'
$ws.Range("J460").Value = 'Statement-non-opinion'
$ws.Range("I485").Value = 'sd'
$ws.Range("J485").Value = 'Statement-non-opinion'
